$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 11 "Marking": Right 5 -> 4, Wrong -1 -> -2
$ws.Range("B11").Value = 4
$ws.Range("C11").Value = -2

# Row 12 "Total": Right 110 -> 88, Wrong -5 -> -10, Max text "110 / 140" -> "78 / 112"
$ws.Range("B12").Value = 88
$ws.Range("C12").Value = -10
$ws.Range("E12").Value = "78 / 112"
